$wb = $excel.ActiveWorkbook

# ALC!2 - Mercury Rising
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 605.4545000000001
$ws.Range("J2").Value = 822.25
$ws.Range("L2").Value = 822.25
$ws.Range("N2").Value = -1048.25

# ALC!28 - The Writing Is Not on the Wall
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1112.7142
$ws.Range("I28").Value = 941.2727
$ws.Range("J28").Value = 1301.3
$ws.Range("K28").Value = 941.2727
$ws.Range("L28").Value = 1301.3
$ws.Range("M28").Value = -456.2727
$ws.Range("N28").Value = -2271.3

# ALC!38 - Just Give Him a Serum
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4622.4443
$ws.Range("I38").Value = 2219.2
$ws.Range("J38").Value = 6339.048
$ws.Range("K38").Value = 6657.599999999999
$ws.Range("L38").Value = 19017.144
$ws.Range("M38").Value = -6285.599999999999
$ws.Range("N38").Value = -19761.144

# ALC!76 - Warding Off Temptation
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9386.385
$ws.Range("I76").Value = 10841.667
$ws.Range("K76").Value = 10841.667
$ws.Range("M76").Value = -10526.667

# ALC!79 - The Garden of Arcane Delights (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 9386.385
$ws.Range("I79").Value = 10841.667
$ws.Range("K79").Value = 10841.667
$ws.Range("M79").Value = -9749.666999999999

# ALC!98 - The Dotted Line
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 960.0714
$ws.Range("I98").Value = 1007
$ws.Range("K98").Value = 1007
$ws.Range("M98").Value = 491

# ALC!121 - Mindful Medicine
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 5449.5
$ws.Range("J121").Value = 5449.5
$ws.Range("L121").Value = 16348.5
$ws.Range("N121").Value = -19842.5

# ALC!122 - Wishful Inking
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 960.0714
$ws.Range("I122").Value = 1007
$ws.Range("K122").Value = 3021
$ws.Range("M122").Value = -571

# ARM!32 - Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3697.757
$ws.Range("I32").Value = 3494.4285
$ws.Range("K32").Value = 3494.4285
$ws.Range("M32").Value = -3207.4285

# ARM!97 - Ore for Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1203.8966
$ws.Range("I97").Value = 1165.875
$ws.Range("J97").Value = 1386.4
$ws.Range("K97").Value = 1165.875
$ws.Range("L97").Value = 1386.4
$ws.Range("M97").Value = -669.875
$ws.Range("N97").Value = -2378.4

# ARM!132 - Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 23835.281
$ws.Range("I132").Value = 1842.3793
$ws.Range("K132").Value = 5527.1379
$ws.Range("M132").Value = -2997.1379

# BSM!20 - Smelt and Dealt
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1823.8125
$ws.Range("I20").Value = 1840.25
$ws.Range("J20").Value = 1774.5
$ws.Range("K20").Value = 1840.25
$ws.Range("L20").Value = 1774.5
$ws.Range("M20").Value = -1593.25
$ws.Range("N20").Value = -2268.5

# BSM!76 - Keep Up with the Mechanics
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 162458.67
$ws.Range("J76").Value = 162458.67
$ws.Range("L76").Value = 162458.67
$ws.Range("N76").Value = -163088.67

# BSM!79 - Unconventional Weaponry (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 162458.67
$ws.Range("J79").Value = 162458.67
$ws.Range("L79").Value = 162458.67
$ws.Range("N79").Value = -164642.67

# BSM!94 - High Steal
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1106.8276
$ws.Range("J94").Value = 1799.6923
$ws.Range("L94").Value = 1799.6923
$ws.Range("N94").Value = -2701.6923

# BSM!107 - The Gold Experience
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2617.5757
$ws.Range("I107").Value = 2363.7334
$ws.Range("K107").Value = 2363.7334
$ws.Range("M107").Value = -443.7334000000001

# CRP!7 - Gridania's Got Talent
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 105.73684
$ws.Range("I7").Value = 39.125
$ws.Range("K7").Value = 39.125
$ws.Range("M7").Value = 73.875

# CUL!63 - The Next to Last Supper
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null

# CUL!66 - Nostalgia through the Stomach (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null

# CUL!86 - Let's Not Get Sappy
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1725.4445
$ws.Range("I86").Value = 1259.75
$ws.Range("J86").Value = 2098
$ws.Range("K86").Value = 3779.25
$ws.Range("L86").Value = 6294
$ws.Range("M86").Value = -2593.25
$ws.Range("N86").Value = -8666

# CUL!89 - Luxury Spillover (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1725.4445
$ws.Range("I89").Value = 1259.75
$ws.Range("J89").Value = 2098
$ws.Range("K89").Value = 11337.75
$ws.Range("L89").Value = 18882
$ws.Range("M89").Value = -5409.75
$ws.Range("N89").Value = -30738

# CUL!104 - Fits to a Tea
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 2899.4
$ws.Range("J104").Value = 2899.4
$ws.Range("L104").Value = 8698.200000000001
$ws.Range("N104").Value = -13940.2

# GSM!80 - Needs More Prayerbell
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3919.8
$ws.Range("J80").Value = 4099.75
$ws.Range("L80").Value = 4099.75
$ws.Range("N80").Value = -6095.75

# GSM!83 - With a Noise That Reaches Heaven (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3919.8
$ws.Range("J83").Value = 4099.75
$ws.Range("L83").Value = 20498.75
$ws.Range("N83").Value = -30482.75

# GSM!113 - Copious Crystal Cannons
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3496.4194
$ws.Range("J113").Value = 4499.4
$ws.Range("L113").Value = 4499.4
$ws.Range("N113").Value = -8839.4

# GSM!123 - Workplace Workout
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 58166.5
$ws.Range("J123").Value = 66333.336
$ws.Range("L123").Value = 66333.336
$ws.Range("N123").Value = -71233.336

# LTW!3 - Underneath It All
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

# LTW!15 - The Bards' Guards
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null

# LTW!22 - Skin off Their Backs
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 955.3043
$ws.Range("I22").Value = 780.35486
$ws.Range("J22").Value = 1316.8667
$ws.Range("K22").Value = 780.35486
$ws.Range("L22").Value = 1316.8667
$ws.Range("M22").Value = -485.35486
$ws.Range("N22").Value = -1906.8667

# LTW!27 - Fire and Hide
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 955.3043
$ws.Range("I27").Value = 780.35486
$ws.Range("J27").Value = 1316.8667
$ws.Range("K27").Value = 780.35486
$ws.Range("L27").Value = 1316.8667
$ws.Range("M27").Value = -673.35486
$ws.Range("N27").Value = -1530.8667

# LTW!46 - Supply Side Logic
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2149.15
$ws.Range("I46").Value = 999.9091
$ws.Range("J46").Value = 3553.7778
$ws.Range("K46").Value = 999.9091
$ws.Range("L46").Value = 3553.7778
$ws.Range("M46").Value = -811.9091
$ws.Range("N46").Value = -3929.7778

# LTW!55 - It's Not a Job, It's a Calling
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 941.36
$ws.Range("I55").Value = 576.6923
$ws.Range("J55").Value = 1336.4166
$ws.Range("K55").Value = 576.6923
$ws.Range("L55").Value = 1336.4166
$ws.Range("M55").Value = -403.6923
$ws.Range("N55").Value = -1682.4166

# LTW!61 - Spelling Me Softly
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3042.6191
$ws.Range("I61").Value = 2520.9473
$ws.Range("J61").Value = 7998.5
$ws.Range("K61").Value = 2520.9473
$ws.Range("L61").Value = 7998.5
$ws.Range("M61").Value = -2318.9473
$ws.Range("N61").Value = -8402.5

# LTW!113 - Peace in Rest
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3042.6191
$ws.Range("I113").Value = 2520.9473
$ws.Range("J113").Value = 7998.5
$ws.Range("K113").Value = 2520.9473
$ws.Range("L113").Value = 7998.5
$ws.Range("M113").Value = -350.9472999999998
$ws.Range("N113").Value = -12338.5

# LTW!122 - Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4123.1875
$ws.Range("I122").Value = 3506
$ws.Range("J122").Value = 5974.75
$ws.Range("K122").Value = 10518
$ws.Range("L122").Value = 17924.25
$ws.Range("M122").Value = -8068
$ws.Range("N122").Value = -22824.25

# LTW!136 - Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4133.75
$ws.Range("I136").Value = 4133.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12401.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9851.25
$ws.Range("N136").Value = $null

# WVR!96 - Skills on Display
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3937.0417
$ws.Range("I96").Value = 3493.7273
$ws.Range("J96").Value = 4312.154
$ws.Range("K96").Value = 3493.7273
$ws.Range("L96").Value = 4312.154
$ws.Range("M96").Value = -2120.7273
$ws.Range("N96").Value = -7058.154

# WVR!107 - Flax Wax
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1315.8572
$ws.Range("I107").Value = 1147.4546
$ws.Range("K107").Value = 3442.3638
$ws.Range("M107").Value = -1522.3638

# WVR!136 - Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4561.6665
$ws.Range("I136").Value = 2607.442
$ws.Range("K136").Value = 7822.326
$ws.Range("M136").Value = -5272.326
